$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.691.15'
$ws.Range('E2').Value = '  +4.68%  '
$ws.Range('D3').Value = '2.298.83'
$ws.Range('E3').Value = '  +2.08%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '304.32'
$ws.Range('E5').Value = '  +1.50%  '
$ws.Range('D6').Value = '101.52'
$ws.Range('E6').Value = '  +9.41%  '
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '0.522'
$ws.Range('E9').Value = '  +3.31%  '
$ws.Range('D10').Value = '36.48'
$ws.Range('E10').Value = '  +7.97%  '
$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').Value = '7.33'
$ws.Range('E12').Value = '  +2.54%  '
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '2.650.16'
$ws.Range('E14').Value = '  +1.87%  '
$ws.Range('D15').Value = '2.301.88'
$ws.Range('E15').Value = '  +1.82%  '
$ws.Range('D16').Value = '13.84'
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').Value = '0.811'
$ws.Range('E17').Value = '  +2.45%  '
$ws.Range('D18').Value = '46.676.32'
$ws.Range('E18').Value = '  +4.55%  '
$ws.Range('D19').Value = '12.96'
$ws.Range('E19').Value = '  +3.34%  '
$ws.Range('D20').Value = '0.0₃0942'
$ws.Range('E20').Value = '  +3.36%  '
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').Value = '65.86'
$ws.Range('E22').Value = '  +1.54%  '
$ws.Range('D23').Value = '249.33'
$ws.Range('E23').Value = '  +4.77%  '
$ws.Range('E24').Value = '  +0.90%  '
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').Value = '1.92'
$ws.Range('E26').Value = '  +1.83%  '
$ws.Range('D27').Value = '42.26'
$ws.Range('E27').Value = '  +7.72%  '
$ws.Range('D28').Value = '2.21'
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('D29').Value = '9.89'
$ws.Range('E29').Value = '  +4.03%  '
$ws.Range('D30').Value = '20.01'
$ws.Range('E30').Value = '  +2.69%  '
$ws.Range('E31').Value = '  +11.04%  '
$ws.Range('D32').Value = '5.59'
$ws.Range('E32').Value = '  +0.91%  '
$ws.Range('D33').Value = '147.52'
$ws.Range('E33').Value = '  -3.58%  '
$ws.Range('D34').Value = '0.0790'
$ws.Range('E34').Value = '  +2.24%  '
$ws.Range('D35').Value = '3.28'
$ws.Range('E35').Value = '  +12.86%  '
$ws.Range('E36').Value = '  +10.12%  '
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = '1.77'
$ws.Range('E38').Value = '  +3.59%  '
$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').Value = '15.96'
$ws.Range('E39').Value = '  +15.73%  '
$ws.Range('D40').Value = '3.97'
$ws.Range('E40').Value = '  +7.58%  '
$ws.Range('D41').Value = '3.35'
$ws.Range('E41').Value = '  +3.37%  '
$ws.Range('D42').Value = '0.0301'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').Value = '1.95'
$ws.Range('E44').Value = '  +7.77%  '
$ws.Range('D45').Value = '1.824.05'
$ws.Range('E45').Value = '  +2.06%  '
$ws.Range('D46').Value = '88.72'
$ws.Range('E46').Value = '  +17.67%  '
$ws.Range('D47').Value = '0.194'
$ws.Range('E47').Value = '  +2.54%  '
$ws.Range('D48').Value = '73.36'
$ws.Range('E48').Value = '  +5.54%  '
$ws.Range('D49').Value = '4.92'
$ws.Range('E49').Value = '  +6.07%  '
$ws.Range('D50').Value = '95.47'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').Value = '2.524.85'
$ws.Range('E51').Value = '  +1.80%  '
